$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells that would otherwise be parsed as numbers
# by Excel (e.g. "1.00", "209.80"), so they stay as text matching the source data.
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '26.390.88'
$ws.Range('E2').Value = '  -1.19%  '
$ws.Range('D3').Value = '1.588.99'
$ws.Range('E3').Value = '  -0.70%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '209.80'
$ws.Range('E5').Value = '  -0.69%  '
$ws.Range('D6').Value = '0.505'
$ws.Range('E6').Value = '  -1.12%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').Value = '0.0613'
$ws.Range('E8').Value = '  -0.88%  '
$ws.Range('D9').Value = '0.246'
$ws.Range('E9').Value = '  -0.36%  '
$ws.Range('D10').Value = '19.54'
$ws.Range('E10').Value = '  -0.89%  '
$ws.Range('E11').Value = '  +0.01%  '
$ws.Range('D12').Value = '1.812.90'
$ws.Range('E12').Value = '  -0.66%  '
$ws.Range('D13').Value = '1.585.46'
$ws.Range('E13').Value = '  -1.24%  '
$ws.Range('D14').Value = '4.03'
$ws.Range('E14').Value = '  -0.12%  '
$ws.Range('D15').Value = '0.518'
$ws.Range('E15').Value = '  -1.03%  '
$ws.Range('D16').Value = '64.32'
$ws.Range('E16').Value = '  -1.31%  '
$ws.Range('D17').Value = '26.386.60'
$ws.Range('E17').Value = '  -1.10%  '
$ws.Range('D18').Value = '0.0₃0740'
$ws.Range('E18').Value = '  -0.80%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').Value = '7.21'
$ws.Range('E19').Value = '  -0.23%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').Value = '1.00'
$ws.Range('E20').Value = '  -0.13%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').Value = '207.07'
$ws.Range('E21').Value = '  -1.48%  '
$ws.Range('D22').Value = '4.29'
$ws.Range('E22').Value = '  -0.25%  '
$ws.Range('D23').Value = '2.22'
$ws.Range('E23').Value = '  -3.29%  '
$ws.Range('D24').Value = '8.83'
$ws.Range('E24').Value = '  -1.28%  '
$ws.Range('D25').Value = '144.68'
$ws.Range('E25').Value = '  +1.22%  '
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').Value = '7.06'
$ws.Range('E27').Value = '  -0.72%  '
$ws.Range('D28').Value = '0.114'
$ws.Range('E28').Value = '  -0.55%  '
$ws.Range('D29').Value = '15.33'
$ws.Range('E29').Value = '  -0.87%  '
$ws.Range('D30').Value = '0.0505'
$ws.Range('E30').Value = '  -1.87%  '
$ws.Range('E31').Value = '  -0.63%  '
$ws.Range('D32').Value = '3.24'
$ws.Range('E32').Value = '  -0.58%  '
$ws.Range('D33').Value = '2.95'
$ws.Range('E33').Value = '  -0.71%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').Value = '1.27'
$ws.Range('E34').Value = '  +15.02%  '
$ws.Range('B35').Value = 'Maker'
$ws.Range('C35').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D35').Value = '1.287.34'
$ws.Range('E35').Value = '  -0.51%  '
$ws.Range('E36').Value = '  +0.53%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '0.599'
$ws.Range('E37').Value = '  -1.11%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').Value = '1.48'
$ws.Range('E38').Value = '  -1.07%  '
$ws.Range('D39').Value = '0.0168'
$ws.Range('E39').Value = '  -1.50%  '
$ws.Range('D40').Value = '0.820'
$ws.Range('E40').Value = '  -0.23%  '
$ws.Range('D41').Value = '5.45'
$ws.Range('E41').Value = '  +0.77%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '0.768'
$ws.Range('E42').Value = '  -1.43%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').Value = '2.14'
$ws.Range('E43').Value = '  -4.28%  '
$ws.Range('D44').Value = '62.62'
$ws.Range('E44').Value = '  -0.83%  '
$ws.Range('D45').Value = '1.726.07'
$ws.Range('E45').Value = '  -0.50%  '
$ws.Range('D46').Value = '89.09'
$ws.Range('E46').Value = '  -2.36%  '
$ws.Range('D47').Value = '1.57'
$ws.Range('E47').Value = '  +0.09%  '
$ws.Range('D48').Value = '0.102'
$ws.Range('E48').Value = '  +1.19%  '
$ws.Range('E49').Value = '  -0.55%  '
$ws.Range('B50').Value = 'USDD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D50').Value = '1.00'
$ws.Range('E50').Value = '  -0.14%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '7.45'
$ws.Range('E51').Value = '  +0.69%  '
